$wb = $excel.ActiveWorkbook

foreach ($idx in 1,4) {
    $ws = $wb.Worksheets.Item($idx)

    $ws.Range("F2").Value = 447

    $ws.Range("A3").Copy()
    $ws.Range("A4").PasteSpecial(-4122)

    $ws.Range("A4").Value = 3
    $ws.Range("B4").NumberFormat = "@"
    $ws.Range("B4").Value = "2024-08-10"
    $ws.Range("C4").Value = "丽水·CCAC动漫七夕（回馈展）"
    $ws.Range("D4").Value = "中东路848号(解放街交汇) 飞达国际大酒店"
    $ws.Range("E4").Value = "2024.08.10 09:00-08.10 17:00"
    $ws.Range("F4").Value = 2
    $ws.Range("G4").Value = 29.9
    $ws.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=86567"
    $ws.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202405/tsOzbBRx1717015539538.png"
}
